# Refresh the cryptocurrency price/volume snapshot (and two rank swaps)
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.204.14"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "2.321.01"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'543.84"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "'130.76"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "2.319.04"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "2.734.05"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "60.160.80"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "2.303.92"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "'10.54"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "'313.51"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "'6.60"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D24").Value = "'63.84"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("E29").Value = "  +5.79%  "
$ws.Range("D30").Value = "'172.05"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'1.72"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "0.0₃0731"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  +7.05%  "
$ws.Range("D35").Value = "'0.379"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'318.20"
$ws.Range("E40").Value = "  +7.59%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'37.89"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.53"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").Value = "'137.48"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").Value = "'3.45"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").Value = "'18.95"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").Value = "'0.564"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0233"
$ws.Range("E48").Value = "  +23.10%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "'0.0492"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  +0.30%  "
